$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point correction on A14 (re-write from source task)
$ws.Range("A14").Value = 45866.62530247685

# Append new row 15 with the latest sensor reading
$ws.Range("A15").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("A15").Value = 45866.66691041316
$ws.Range("B15").Value = 2025
$ws.Range("C15").Value = 31
$ws.Range("D15").Value = 20.28
$ws.Range("E15").Value = 71.92
$ws.Range("F15").Value = 325.2
$ws.Range("G15").Value = 12.45
$ws.Range("H15").Value = "ESE"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "16:00:21"
